$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns hold plain text in this sheet
# (e.g. "301.25", "-0.55%"), not real numbers/percentages. Excel normally
# auto-converts numeric- and percent-looking text on assignment, so for each
# touched cell we flip the range to Text format first, write the new string,
# then restore the "Normal" style so no stray formatting is left behind.

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '301.25'
$ws.Range("E2").Value = '-0.55%'
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '37.63'
$ws.Range("E3").Value = '8.09%'
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = '5.008'
$ws.Range("E4").Value = '-2.87%'
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '0.07847'
$ws.Range("E5").Value = '1.24%'
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '2.193'
$ws.Range("E6").Value = '-7.82%'
$rng.Style = "Normal"

$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$ws.Range("E7").Value = '-0.07%'
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '4.000'
$ws.Range("E8").Value = '1.47%'
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '0.9109'
$ws.Range("E9").Value = '-2.00%'
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.1869'
$ws.Range("E10").Value = '3.97%'
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '0.09218'
$ws.Range("E11").Value = '-6.53%'
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '0.08451'
$ws.Range("E12").Value = '-1.84%'
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '0.03537'
$ws.Range("E13").Value = '6.47%'
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '0.09937'
$ws.Range("E14").Value = '0.51%'
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '0.001470'
$ws.Range("E15").Value = '-1.75%'
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '0.005645'
$ws.Range("E16").Value = '-1.68%'
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '3.476'
$ws.Range("E17").Value = '0.47%'
$rng.Style = "Normal"

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = '-1.81%'
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = '2.82%'
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '0.1316'
$ws.Range("E20").Value = '-1.35%'
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '4.534'
$ws.Range("E21").Value = '4.87%'
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '0.2243'
$ws.Range("E22").Value = '-2.45%'
$rng.Style = "Normal"

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$ws.Range("E23").Value = '1.33%'
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '0.001229'
$ws.Range("E24").Value = '0.99%'
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '0.004442'
$ws.Range("E25").Value = '-0.42%'
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = '0.0001297'
$ws.Range("E26").Value = '-0.18%'
$rng.Style = "Normal"

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '0.0004747'
$ws.Range("E27").Value = '39.90%'
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = '0.01747'
$ws.Range("E39").Value = '-2.64%'
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = '0.04709'
$ws.Range("E40").Value = '-1.67%'
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = '0.007866'
$ws.Range("E41").Value = '1.58%'
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = '0.1387'
$ws.Range("E42").Value = '-1.61%'
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = '0.007665'
$ws.Range("E43").Value = '7.98%'
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = '0.002295'
$ws.Range("E44").Value = '9.40%'
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '0.01013'
$ws.Range("E45").Value = '10.35%'
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = '0.00006057'
$ws.Range("E46").Value = '-1.01%'
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").Value = '-0.18%'
$rng.Style = "Normal"

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = '8.667'
$ws.Range("E48").Value = '182.98%'
$rng.Style = "Normal"

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = '34.95%'
$rng.Style = "Normal"

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '0.00002096'
$ws.Range("E50").Value = '-0.18%'
$rng.Style = "Normal"

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = '0.0001996'
$ws.Range("E51").Value = '-0.18%'
$rng.Style = "Normal"
